# Leave Card update (1/5/2024 4:46 PM) — add a new period-end row, shift the
# 2024 month-end marker dates by one day (now last-day-of-month instead of
# first-day-of-next-month), post two 1.25 SL/VL accruals for Nov/Dec 2023,
# and record a new SL(1-0-0) leave entry for January 2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# --- 1. Grow the table by one row (A8:K103 -> A8:K104) --------------------
$lo.Resize($ws.Range("A8:K104"))

# Duplicate row 103's cell formatting/formula down into the brand-new row 104
$ws.Range("A103:K103").Copy($ws.Range("A104:K104"))
$ws.Range("A104").Value2 = 45626
$ws.Range("G104").Formula = '=IF(ISBLANK([@EARNED]),"",[@EARNED])'

# --- 2. Post the two pending 1.25 accruals (Nov/Dec 2023) ------------------
$ws.Range("C91").Value2 = 1.25
$ws.Range("C92").Value2 = 1.25

# --- 3. Turn the "01/01/2024" row-marker into a plain "2024" year label ----
# (mirrors the existing 2020/2021/2022/2023 labels elsewhere in the sheet)
$ws.Range("A74").Copy()
$ws.Range("A93").PasteSpecial(-4122)
$ws.Range("A93").Value = "'2024"

# --- 4. New SL(1-0-0) leave entry charged against January 2024 -------------
$ws.Range("B94").Value2 = "SL(1-0-0)"
$ws.Range("H94").Value2 = 1
$ws.Range("K90").Copy()
$ws.Range("K94").PasteSpecial(-4122)
$ws.Range("K94").Value2 = 44930

# --- 5. Shift the 2024 month-end marker dates back by one day --------------
# (these had been the 1st of the following month; now the last day of the
# current month, matching every later row already in the sheet)
$ws.Range("A94").Value2 = 45322
$ws.Range("A95").Value2 = 45351
$ws.Range("A96").Value2 = 45382
$ws.Range("A97").Value2 = 45412
$ws.Range("A98").Value2 = 45443
$ws.Range("A99").Value2 = 45473
$ws.Range("A100").Value2 = 45504
$ws.Range("A101").Value2 = 45535
$ws.Range("A102").Value2 = 45565
$ws.Range("A103").Value2 = 45596

$wb.Application.Calculate()
